$d = $word.ActiveDocument

# Collapse the "<id>...</id>" spans that are currently split across
# three runs (open-tag run, plain-text id run, close-tag run) back
# into a single run, so the text reads as one contiguous token with
# the tag's own (Courier New / 7f6000) character formatting.
#
# Word's Find/Replace collapses a multi-run match into a single run
# using the formatting of the first run in the match, which is
# exactly the behavior the diff shows (the "<id>" run's rPr survives,
# the middle plain-text run's rPr and the "</id>" run's rPr are gone).

$d.Content.Find.Execute("<id>p035v_2</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p035v_2</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p036r_1</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p036r_1</id>", 2) | Out-Null
